# Actualización automática del tracker
# - Rellena el resultado/profit de la fila 60 (Gonzalo Villanueva vs Facundo Juarez)
# - Agrega una nueva fila (67) con el próximo partido a trackear

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fila 60: se conoce el resultado -> "Fallo" / profit -1 ---
$ws.Range("G60").Value = "Fallo"
$ws.Range("H60").Value = -1

# --- Fila 67: nuevo partido agregado al tracker ---
$ws.Range("A67").Value = 14740886

# B67 es una fecha en formato texto ("2025-09-23"); el apóstrofo fuerza texto
# para que Excel no lo interprete como un valor de fecha, y luego se limpia
# el formato para no dejar la celda marcada con "quote prefix".
$ws.Range("B67").Value = "'2025-09-23"
$ws.Range("B67").ClearFormats()

$ws.Range("C67").Value = "Juan Manuel Cerundolo"
$ws.Range("D67").Value = "Sho Shimabukuro"
$ws.Range("E67").Value = "Gana Juan Manuel Cerundolo"
$ws.Range("F67").Value = 1.73

# resultado/profit todavía no se conocen -> quedan como texto vacío
$ws.Range("G67").Value = "'"
$ws.Range("G67").ClearFormats()
$ws.Range("H67").Value = "'"
$ws.Range("H67").ClearFormats()
